# Update countries & provincias Spain
# Applies the 19-May-2020 15:05 data refresh to the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: move country-name labels around.
# Several countries changed rank (their row position), which in the
# underlying workbook model means the *text* shown in a given row changes
# while the row number itself stays put. Several of these are rotations
# among a handful of rows (e.g. Tayikistan jumps to the top of its
# neighbourhood, pushing Cuba / Macedonia / Islandia / Estonia down by one).
#
# Doing this as a naive sequence of assignments risks a transient moment
# where two cells hold identical text (e.g. writing "Estonia" into A89
# while A88 still holds "Estonia"), which would collapse the two distinct
# strings into one. To avoid that, every affected cell is first given a
# unique placeholder, and only once all placeholders are in place do we
# write the real final text.

$placeholders = @{
    "A85"  = "__PH_A85__"
    "A86"  = "__PH_A86__"
    "A87"  = "__PH_A87__"
    "A88"  = "__PH_A88__"
    "A89"  = "__PH_A89__"
    "A107" = "__PH_A107__"
    "A108" = "__PH_A108__"
    "A109" = "__PH_A109__"
    "A110" = "__PH_A110__"
    "A196" = "__PH_A196__"
    "A197" = "__PH_A197__"
    "A209" = "__PH_A209__"
    "A210" = "__PH_A210__"
    "A211" = "__PH_A211__"
    "A215" = "__PH_A215__"
    "A216" = "__PH_A216__"
}

foreach ($addr in $placeholders.Keys) {
    $ws.Range($addr).Value = $placeholders[$addr]
}

$finalNames = @{
    "A85"  = "Tayikistan"
    "A86"  = "Cuba"
    "A87"  = "Republica de Macedonia"
    "A88"  = "Islandia"
    "A89"  = "Estonia"
    "A107" = "Kenia"
    "A108" = "Libano"
    "A109" = "Albania"
    "A110" = "Republica de Chipre"
    "A196" = "Nueva Caledonia"
    "A197" = "Belice"
    "A209" = "Seychelles"
    "A210" = "Groenlandia"
    "A211" = "Montserrat"
    "A215" = "San Bartolome"
    "A216" = "Bonaire, San Eustaquio y Saba"
}

foreach ($addr in $finalNames.Keys) {
    $ws.Range($addr).Value = $finalNames[$addr]
}

# --- Step 2: refresh the case-count statistics (columns B..H) for every
# row whose figures moved along with (or independently of) the re-ranking
# above.

$stats = @{
    4   = @(1550699, 405, 358906, 1099757, 0, 55, 92036)
    14  = @(102287, 1959, 39658, 59460, 0, 13, 3169)
    18  = @(59854, 2509, 31634, 27891, 0, 9, 329)
    49  = @(10733, 34, 4904, 5595, 0, 3, 234)
    51  = @(8604, 18, 5687, 2616, 0, 4, 301)
    61  = @(6399, 19, 5000, 1098, 0, 1, 301)
    85  = @(1936, 207, 0, 1895, 0, 0, 41)
    86  = @(1881, 0, 1505, 297, 0, 0, 79)
    87  = @(1839, 22, 1351, 382, 0, 2, 106)
    88  = @(1802, 0, 1786, 6, 0, 0, 10)
    89  = @(1791, 7, 938, 789, 0, 0, 64)
    101 = @(1143, 37, 80, 1059, 0, 0, 4)
    107 = @(963, 51, 358, 555, 0, 0, 50)
    108 = @(954, 23, 251, 677, 0, 0, 26)
    109 = @(949, 1, 742, 176, 0, 0, 31)
    110 = @(917, 0, 515, 385, 0, 0, 17)
    196 = @(18, 0, 18, 0, 0, 0, 0)
    197 = @(18, 0, 16, 0, 0, 0, 2)
    210 = @(11, 0, 11, 0, 0, 0, 0)
    211 = @(11, 0, 10, 0, 0, 0, 1)
}

foreach ($row in $stats.Keys) {
    $vals = $stats[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
    $ws.Cells.Item($row, 8).Value = $vals[6]
}

# --- Step 3: bump the "last updated" timestamp banner.
$ws.Range("A1").Value = "Datos actualizados a 19 de Mayo de 2020 a las 15:05"
